# Insert a new record row at row 2 of the Initial_Screening sheet
# (ScreenName/SectionName/SubSectionName/FieldName/FieldType/ElementType/ID),
# pushing the existing rows 2-11 down to 3-12, and move the active
# selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2:11 -> 3:12).
$ws.Rows.Item(2).Insert()

# The inserted row copies the header row's bold/fill formatting by default;
# reset A2:G2 back to the plain "Normal" style used elsewhere in the sheet,
# with wrap-text enabled (matches the workbook's existing wrap-text style).
$newRow = $ws.Range("A2:G2")
$newRow.Style = "Normal"
$newRow.WrapText = $true

# Columns H2:I2 should stay empty/unused for this row (no residual style).
$ws.Range("H2:I2").Clear()

# Populate the new record.
$ws.Range("A2").Value = "Initial_Screening"
$ws.Range("B2").Value = "Screening"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "NewBtn"
$ws.Range("E2").Value = "Div"
$ws.Range("F2").Value = "Link"
$ws.Range("G2").Value = "NA"

# Move the selection, matching the saved workbook's cursor position.
$ws.Range("D6").Select() | Out-Null
